# Applies the 2025-12-15 03:30 "All-in-one: update products & pages" edit:
# appends 16 new product rows (rows 51-66) to Sheet1, mirroring the
# automated product-feed pipeline's latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 16 new product rows (50-65, 1-indexed product #, sheet rows 51-66)
$newRows = @(
    @(50, '휴고스 미니 컴팩트 PTC 전기온풍기 소형 히터 캠핑 사무실 발난로', 'https://naver.me/xbjdAyGm', 'https://shop-phinf.pstatic.net/20250926_101/17588598135561I8SR_JPEG/92992735688983672_1838852221.jpg', '미니 전기온풍기 발난로', '작고 강력한 온풍기로 사무실과 캠핑장에서 따뜻함을 빠르게 느껴보세요.', '휴고스 미니 컴팩트 PTC 전기온풍기 추천', '작고 강력한 휴고스 미니 전기온풍기로 사무실과 캠핑 등 다양한 공간을 빠르게 따뜻하게 해 줍니다.', '휴고스,미니온풍기,컴팩트히터,PTC히터,전기온풍기,소형히터,발난로,캠핑히터,사무실히터', '가전/히터'),
    @(51, '르젠 감성 동파방지 무타공 1초발열 절전형 순간발열 PTC 욕실용히터 온풍기 LZH-B770D', 'https://naver.me/5eZBUv8s', 'https://shop-phinf.pstatic.net/20251126_137/1764118554399jnJnU_JPEG/25455801497886398_1776636856.jpg', '절전형 욕실용 순간발열 히터', '추운 욕실도 빠르게 따뜻하게, 전기세 걱정 없이 절전형 발열 기술로 쾌적함을 유지하세요.', '르젠 감성 동파방지 욕실용 히터 LZH-B770D 절전형 구매', '추운 욕실을 빠르게 따뜻하게 만드는 절전형 PTC 순간발열 욕실용 온풍기로 실용적이고 경제적입니다.', '르젠,동파방지,욕실히터,절전형,순간발열,PTC히터,온풍기,욕실온풍기', '가전,욕실용품'),
    @(52, '쿠오레 불멍 전기 가습 히터 가정용 거실 사무실 온풍기 벽난로 인테리어 CH-LH20B', 'https://naver.me/55PhraLO', 'https://shop-phinf.pstatic.net/20251205_157/1764921072519h2VVN_PNG/99053885603308557_650950176.png', '거실용 전기 가습 히터', '추운 겨울, 따뜻한 온기와 촉촉한 습기로 쾌적한 공간을 만들어 보세요.', '쿠오레 불멍 전기 가습 히터 CH-LH20B 온풍기 추천', '쿠오레 불멍 가습 히터는 거실과 사무실에 따뜻한 온기와 촉촉한 습기를 제공합니다. 겨울철 쾌적한 실내 환경에 적합합니다.', '쿠오레,전기히터,가습히터,온풍기,벽난로,거실히터,사무실히터,겨울가전', '가전'),
    @(53, '가정용 욕실 온풍기 화장실 벽걸이 난방기 PTC 히터 2255WF', 'https://naver.me/xQiXMG1Y', 'https://shop-phinf.pstatic.net/20251010_8/17600522426622O3Ko_JPEG/45896581643329479_1789907410.jpg', '욕실 벽걸이 온풍기', '추운 욕실도 걱정 끝, 빠르게 따뜻한 온기로 쾌적한 공간을 만들어드립니다.', '가정용 욕실 온풍기 PTC 히터 2255WF 추천 및 할인', '추운 욕실에 적합한 벽걸이형 온풍기로 빠르게 따뜻한 공간을 제공합니다. 쾌적한 욕실 환경을 원하는 가정에 추천합니다.', '욕실온풍기,가정용온풍기,벽걸이난방기,PTC히터,화장실난방,욕실난방기,온풍기추천', '가전/난방기'),
    @(54, '가정용온풍기 저소음 PTC히터 스탠드 난방기 사무실온풍기 2550', 'https://naver.me/5L7ZHXqX', 'https://shop-phinf.pstatic.net/20251103_37/1762124882953bmRwM_JPEG/78350197805646922_647395038.jpg', '저소음 가정용 온풍기', '추운 날씨에도 조용하고 빠른 온기로 따뜻함을 느껴보세요. 사무실과 가정 어디서나 편안한 난방을 제공합니다.', '가정용온풍기 2550 저소음 PTC히터 난방기 추천', '가정과 사무실에서 사용할 수 있는 저소음 PTC 스탠드 온풍기로 빠르게 따뜻함을 제공합니다.', '가정용온풍기,저소음히터,PTC히터,스탠드난방기,사무실온풍기,효율난방기', '가전/난방기'),
    @(55, '이나프 무선 청소기 MV1 차량용 미니 핸디 소형 원룸 청소기', 'https://naver.me/5UTCEAyE', 'https://shop-phinf.pstatic.net/20251202_14/1764633907026ALdah_JPEG/2418911230288605_1249188675.jpg', '차량용 초경량 무선청소기', '좁은 공간도 깔끔하게! 차량과 원룸 청소 고민을 한 번에 해결하세요.', '이나프 무선 차량용 미니 청소기 MV1 추천', '좁은 공간과 차량, 원룸 청소에 최적화된 무선 핸디 청소기로 간편하게 청소하세요.', '이나프,무선청소기,차량용청소기,미니청소기,핸디청소기,원룸청소기,소형청소기', '가전/청소기'),
    @(56, '드리미 V10s 무선청소기 핸디스틱형 진공청소기 140AW 3자기 흡입 모드 국내 A/S', 'https://naver.me/xhlqzLoH', 'https://shop-phinf.pstatic.net/20250812_259/17549936621890965E_PNG/81486967999075804_618541934.png', '강력 무선 진공청소기', '가벼운 무선으로 집안 구석구석 깔끔하게! 흡입력 걱정 없이 편리한 청소를 경험하세요.', '드리미 V10s 무선청소기 핸디 스틱형 추천 및 할인', '가벼운 무선 청소기로 집안 구석구석 손쉽게 청소하세요. 140AW 강력 흡입력에 3가지 자기 흡입 모드로 효율적입니다.', '드리미,무선청소기,핸디스틱,진공청소기,강력흡입,가벼운청소기,국내용A/S,140AW', '청소기/가전'),
    @(57, '몬스타기어 클래식104 멤브레인 보글보글 사무용 저소음 키보드 투명 실리콘 키스킨 세트', 'https://naver.me/FOZj95EJ', 'https://shop-phinf.pstatic.net/20240215_292/17079859121505xnDj_JPEG/109121739972719035_130031986.jpg', '저소음 사무용 멤브레인 키보드', '조용한 업무 환경을 원한다면 저소음 멤브레인 키보드가 제격입니다. 투명 실리콘 키스킨으로 깨끗함까지 챙겨보세요.', '몬스타기어 클래식104 저소음 멤브레인 키보드 추천', '사무용으로 적합한 저소음 멤브레인 키보드와 투명 실리콘 키스킨 세트로 조용하고 깔끔한 업무 환경을 만들어줍니다.', '몬스타기어,저소음키보드,멤브레인키보드,사무용키보드,투명키스킨,키스킨세트,실리콘키스킨,조용한키보드', '사무용,키보드'),
    @(58, '가츠 클래식104 보글보글 5세대 멤브레인 저소음 풀윤활 게이밍 사무용 유선 키보드', 'https://naver.me/58jUNQIy', 'https://shop-phinf.pstatic.net/20240215_292/17079859121505xnDj_JPEG/109121739972719035_130031986.jpg', '저소음 멤브레인 게이밍 키보드', '소음 걱정 없이 게임과 업무에 집중하세요. 부드러운 타건감과 풀윤활 기능이 손끝부터 편안함을 선사합니다.', '가츠 클래식104 멤브레인 저소음 키보드 추천', '저소음 멤브레인과 풀윤활 기능으로 부드럽고 편안한 타건감을 제공하는 게이밍 및 사무용 유선 키보드입니다.', '가츠,클래식104,멤브레인 키보드,저소음 키보드,게이밍 키보드,사무용 키보드,유선 키보드,풀윤활 키보드', '키보드/컴퓨터'),
    @(59, '몬스타기어 X108 보글보글 5세대 멤브레인 사무용 저소음 풀윤활 키보드 리퍼 리퍼비시', 'https://naver.me/xQiveUCS', 'https://shop-phinf.pstatic.net/20241129_297/17328483143964YXbe_JPEG/14921095412492371_319647656.jpg', '사무용 저소음 멤브레인 키보드', '조용한 타건감으로 업무 집중을 도와주는 5세대 멤브레인 키보드입니다.', '몬스타기어 X108 보글보글 5세대 멤브레인 키보드 추천', '저소음 타입의 5세대 멤브레인 키보드로 사무용에 적합하며, 조용한 타건감이 업무 집중을 지원합니다.', '몬스타기어, 멤브레인 키보드, 사무용 키보드, 저소음 키보드, 풀윤활, 리퍼비시 키보드, X108', '사무용/키보드'),
    @(60, '호무로 그래피트 알러지프리 극세사 자가발열 겨울 패드 SS', 'https://naver.me/xnrHIZY3', 'https://shop-phinf.pstatic.net/20241217_285/1734421817729xsGjl_JPEG/68554763848718013_580686861.jpg', '알러지 방지 자가발열 패드', '민감한 피부도 걱정 없는 알러지프리 자가발열 패드로 따뜻한 겨울을 경험해보세요.', '호무로 그래피트 알러지프리 자가발열 겨울 패드 SS 추천', '민감한 피부도 안심하고 사용할 수 있는 알러지프리 자가발열 극세사 패드로 따뜻하고 편안한 겨울을 보내세요.', '호무로,그래피트,알러지프리,자가발열,겨울패드,극세사,민감피부,따뜻한침구', '침구/겨울'),
    @(61, '호이데코 코지 단모 러그 [워셔블] - 50x100', 'https://naver.me/Fhf1qzbE', 'https://shop-phinf.pstatic.net/20250324_277/1742798947884NpXGU_JPEG/76931733998516137_1175007394.jpg', '포근한 코지 단모 러그', '차가운 바닥 걱정 없이 따뜻한 공간을 만들어줄 워셔블 러그입니다.', '호이데코 코지 단모 러그 50x100 워셔블 러그 추천', '차가운 바닥을 따뜻하게 보호하는 워셔블 러그로 실내 공간을 아늑하게 만들어줍니다. 세탁이 쉬워 청결하게 사용할 수 있습니다.', '러그,단모러그,워셔블러그,따뜻한러그,호이데코,실내용러그,러그추천', '인테리어/러그'),
    @(62, '원목 침대 프레임 매트리스 깔판 저상형 침대 받침대 소나무 2단 접이식 깔판 미니싱글', 'https://naver.me/5CF3SNoV', 'https://shop-phinf.pstatic.net/20220215_148/1644907248123CQ5MB_JPEG/46043136798435759_475255463.jpg', '저상형 원목 침대 프레임', '견고한 소나무 2단 접이식 디자인으로 공간 활용을 극대화하고 쾌적한 잠자리를 제공합니다.', '원목 소나무 2단 접이식 저상형 침대 프레임 미니싱글 추천', '견고한 소나무 소재의 2단 접이식 원목 침대 프레임으로 공간 활용과 쾌적한 잠자리를 모두 만족시키는 미니싱글 침대 받침대입니다.', '원목침대,소나무침대,저상형침대,접이식침대,미니싱글침대,침대프레임,침대받침대', '가구/침대'),
    @(63, '모노 러그 거실 카페트 먼지없는 사계절 단모 대형 온감 워셔블 아트 100x150', 'https://naver.me/FMT4y2f2', 'https://shop-phinf.pstatic.net/20250103_186/1735853909606RTMwG_JPEG/16485232410076329_319247054.jpg', '사계절 먼지없는 대형 러그', '먼지 걱정 없이 쾌적한 거실을 만들어주는 사계절용 대형 러그로 집안을 따뜻하고 아늑하게 꾸며보세요.', '모노 러그 거실용 대형 사계절 먼지없는 단모 카페트 추천', '먼지 없는 사계절용 대형 러그로 쾌적하고 따뜻한 거실 인테리어에 적합하며 세탁도 간편한 워셔블 제품입니다.', '모노러그,거실러그,사계절러그,대형카페트,먼지없는러그,단모러그,워셔블러그,온감러그,아트러그', '거실/인테리어'),
    @(64, '사계절 이불 포근한 차렵 겨울이불 알러지케어 에코 옥수수솜 단품 S', 'https://naver.me/FK0uYd8y', 'https://shop-phinf.pstatic.net/20240911_232/1726043960560guMeS_JPEG/29600834309618471_1085720924.jpg', '사계절 포근 차렵이불', '알러지 걱정 없이 사계절 내내 포근하게 감싸주는 옥수수솜 이불로 쾌적한 잠자리를 완성하세요.', '사계절 이불 포근한 차렵 겨울이불 알러지케어 추천', '알러지 걱정 없는 에코 옥수수솜으로 사계절 내내 포근하게 사용할 수 있는 차렵 이불입니다.', '사계절이불,차렵이불,알러지케어,옥수수솜이불,포근한이불,겨울이불,에코이불,단품이불,S사이즈', '침구/이불'),
    @(65, '핸드폰 사진인화 액자 사진인쇄 사진출력 웨딩 제작 만들기 미니 9x9cm', 'https://naver.me/FjCSRZ3S', 'https://shop-phinf.pstatic.net/20241114_112/1731563190364xHlO5_JPEG/59945680989044442_289981561.jpg', '미니 사진인화 액자', '소중한 추억을 작고 귀여운 액자에! 특별한 순간을 간편하게 간직해보세요.', '핸드폰 사진인화 미니 액자 9x9cm 제작 추천', '작고 귀여운 9x9cm 미니 액자에 핸드폰 사진을 인화하여 웨딩 등 특별한 순간을 간편하게 간직할 수 있습니다.', '핸드폰 사진인화, 미니 액자, 사진출력, 웨딩 사진, 사진인쇄, 소형 액자, 추억 보관', '사진인화/액자')
)

$startRow = 51
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

# Restore the selection to A1 (matches the refreshed sheet view after the
# pipeline regenerates/reopens the workbook).
[void]$ws.Range("A1").Select()

Write-Output ("Added " + $newRows.Count + " product rows; sheet now spans " + $ws.UsedRange.Address())
